$d = $word.ActiveDocument

$d.Content.Find.Execute("26 uL ddH2O", $true, $false, $false, $false, $false,
                         $true, 1, $false, "33 uL ddH2O", 2)

$d.Content.Find.Execute("4 uL 10x T4 DNA Ligase Buffer", $true, $false, $false, $false, $false,
                         $true, 1, $false, "5.5 uL 10x T4 DNA Ligase Buffer", 2)

$d.Content.Find.Execute("2 uL frag2 (back1)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "5.5 uL frag2 (back1)", 2)

$d.Content.Find.Execute("2 uL BsaI", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2.75 uL BsaI", 2)
